$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price (D) and Volume (E) columns hold text-formatted values (e.g. "29.182.99",
# "  +0.90%  "), matching the source data feed which stores these as text/inline strings,
# not native numbers. Force column D cells to Text format before assignment so Excel
# does not auto-convert numeric-looking strings (e.g. "7.694", "9.780") into numbers,
# which would silently strip significant trailing zeros and change the cell type.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.182.99'
$ws.Range('E2').Value = '  +0.90%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.940.36'
$ws.Range('E3').Value = '  +2.55%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9976'
$ws.Range('E4').Value = '  -0.47%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '326.44'
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9979'
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4612'
$ws.Range('E7').Value = '  +0.59%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3904'
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07872'
$ws.Range('E9').Value = '  +0.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9967'
$ws.Range('E10').Value = '  +0.87%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '22.31'
$ws.Range('E11').Value = '  +1.86%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.923.18'
$ws.Range('E12').Value = '  +2.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.844'
$ws.Range('E13').Value = '  +2.55%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.109'
$ws.Range('E14').Value = '  +1.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07051'
$ws.Range('E15').Value = '  +1.50%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '87.72'
$ws.Range('E16').Value = '  -0.31%  '
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009949'
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.13'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9992'
$ws.Range('E20').Value = '  -0.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '29.228.06'
$ws.Range('E21').Value = '  +1.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.520'
$ws.Range('E22').Value = '  +4.41%  '
$ws.Range('E23').Value = '  +1.82%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.165.93'
$ws.Range('E24').Value = '  +2.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.092'
$ws.Range('E25').Value = '  +1.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '156.14'
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('E27').Value = '  +0.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.904'
$ws.Range('E28').Value = '  -0.49%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '118.77'
$ws.Range('E29').Value = '  +1.12%  '
$ws.Range('E30').Value = '  -2.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09329'
$ws.Range('E31').Value = '  -0.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.8951'
$ws.Range('E32').Value = '  -1.67%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.217'
$ws.Range('E33').Value = '  -1.48%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.323'
$ws.Range('E34').Value = '  -0.80%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.148'
$ws.Range('E35').Value = '  -4.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.000003737'
$ws.Range('E36').Value = '  +135.27%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05794'
$ws.Range('E37').Value = '  +0.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.172'
$ws.Range('E38').Value = '  -1.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02106'
$ws.Range('E39').Value = '  +1.66%  '
$ws.Range('B40').Value = 'Frax'
$ws.Range('C40').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9979'
$ws.Range('E40').Value = '  -0.34%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.694'
$ws.Range('E41').Value = '  -0.61%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5712'
$ws.Range('E42').Value = '  +0.53%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1816'
$ws.Range('E43').Value = '  +2.59%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '9.780'
$ws.Range('E44').Value = '  +0.18%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '11.98'
$ws.Range('E45').Value = '  +0.91%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.210'
$ws.Range('E46').Value = '  -1.67%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5335'
$ws.Range('E47').Value = '  -0.28%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.06936'
$ws.Range('E48').Value = '  -1.49%  '
$ws.Range('B49').Value = 'MXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.589'
$ws.Range('E49').Value = '  +2.41%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.850'
$ws.Range('E50').Value = '  +0.59%  '
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '113.34'
$ws.Range('E51').Value = '  +0.82%  '
